$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the n-channel MOSFET row (row 10): BSH103,235 (SOT-23-3, C88711)
# is replaced by the new part 2N7002K-7 (SOT-23-3, C705095), with
# designators Q1-Q4 renumbered to Q5-Q8.
$ws.Range("A10").Value = "2N7002K-7"
$ws.Range("B10").Value = "Q5, Q6, Q7, Q8"
$ws.Range("D10").Value = "C705095"

# Move the active selection to C11, matching the saved cursor position.
$ws.Range("C11").Select()
